$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (78) down into the two
# new rows so the date column keeps its date style (s="1"), then overwrite
# the values.
$ws.Range("A78:G78").Copy()
$ws.Range("A79:G79").PasteSpecial(-4122)
$ws.Range("A78:G78").Copy()
$ws.Range("A80:G80").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 79 — 2025-06-01
$ws.Range("A79").Value = 45809
$ws.Range("B79").Value = -0.436
$ws.Range("C79").Value = -0.223
$ws.Range("D79").Value = 0.172
$ws.Range("E79").Value = 0.293
$ws.Range("F79").Value = 0.072
$ws.Range("G79").Value = 79.67

# Row 80 — 2025-07-01
$ws.Range("A80").Value = 45839
$ws.Range("B80").Value = -0.54
$ws.Range("C80").Value = -0.258
$ws.Range("D80").Value = 0.162
$ws.Range("E80").Value = 0.272
$ws.Range("F80").Value = 0.193
$ws.Range("G80").Value = 79.59
